$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182),
    @(3.230985683306322, 1.667794583268128, 3.900430680208489, 8.660232485948974, 17.45944343273191),
    @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759),
    @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797),
    @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 31.61296591696135),
    @(1.459612070389937, 0.04240448674262143, 3.900430680208489, 8.660232485948974, 14.06267972329002)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
